$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data rows (rows 2-13), columns A-H:
# A: Instance, B: Method, C: Wmax, D: Wmin, E: Wmax-Wmin, F: Total W, G: Average W, H: Execution Time
$data = @(
    @("40_homogeneous",  "deterministic", 508.72, 500.4,   8.32,   1009.12, 504.56,            0.003000736236572266),
    @("40_homogeneous",  "randomized",    514.89, 513.13,  1.75,   1028.02, 514.01,             1.670448303222656),
    @("40_heterogeneous","deterministic", 641.67, 450.95,  190.72, 1092.62, 546.3099999999999,  0.002999305725097656),
    @("40_heterogeneous","randomized",    641.03, 480.89,  160.14, 1121.93, 560.96,             1.758956670761108),
    @("60_homogeneous",  "deterministic", 546.88, 532.5700000000001, 14.32, 1622.81, 540.9400000000001, 0.01012277603149414),
    @("60_homogeneous",  "randomized",    556.29, 547.04,  9.25,   1657.63, 552.54,             3.73978066444397),
    @("60_heterogeneous","deterministic", 631.05, 491.1,   139.95, 1746.6,  582.2,              0.01000070571899414),
    @("60_heterogeneous","randomized",    652.01, 509.49,  142.52, 1807.74, 602.58,             4.066461801528931),
    @("80_homogeneous",  "deterministic", 532.41, 514.28,  18.13,  2086.54, 521.63,             0.01059842109680176),
    @("80_homogeneous",  "randomized",    541,    536.12,  4.88,   2156.28, 539.0700000000001,  6.165876150131226),
    @("80_heterogeneous","deterministic", 670.96, 456.45,  214.51, 2256.02, 564.01,             0.01103901863098145),
    @("80_heterogeneous","randomized",    693.26, 442.47,  250.79, 2346.15, 586.54,             6.475290536880493)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row = $row + 1
}
